$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.992.44'
$ws.Range('E2').Value = '  +1.38%  '
$ws.Range('D3').Value = '1.847.24'
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.013'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '309.95'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4774'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3677'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('E9').Value = '  +1.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9309'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.74'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07728'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.336'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.763.31'
$ws.Range('E14').Value = '  -3.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.440'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.87'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('D20').Value = '27.030.29'
$ws.Range('E20').Value = '  +1.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.51'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.061'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('E23').Value = '  +0.89%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.931'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.21'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.003'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.36'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.009'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08903'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.295'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.175'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7460'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.505'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.734'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.114'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.10%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01955'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.29%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05275'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.980'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.18%  '
$ws.Range('E40').Value = '  +3.68%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.019'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1514'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.223'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4753'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.00%  '
$ws.Range('E46').Value = '  +0.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.19'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.28%  '
$ws.Range('E48').Value = '  +2.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '66.16'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06063'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8878'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.07%  '
